$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 20241008
# B2 loses its formatting (becomes default/Normal style) once the value is entered
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = 2218.31336
$ws.Range("C2").Value = 2224.47
$ws.Range("D2").Formula = "=100*(B2-C2)/C2"
$ws.Range("E2").Value = 180
$ws.Range("F2").Value = "No open date"

# --- Row 3 ---
$ws.Range("A3").Value = 20241009
$ws.Range("C3").Value = 2224.47
$ws.Range("E3").Value = 180
$ws.Range("F3").Value = "No open date"

# --- Row 4 ---
$ws.Range("A4").Value = 20241010
$ws.Range("C4").Value = 2224.47
$ws.Range("E4").Value = 180
$ws.Range("F4").Value = 20241010

# --- Row 7: C7:F7 are fully cleared back to blank/unstyled cells ---
$ws.Range("C7:F7").Style = "Normal"
$ws.Range("C7:F7").ClearContents()

# --- Selection moves to E7 ---
$ws.Range("E7").Select()
